# Generate Report for Handoff
# Update the "47407ed6-37c6-4239-8f04-d64dee116c62.md" row (row 3) across the
# Overview, zh-cn and de-de sheets: status flips from "Handed back: in sync
# with en-US" back to "In Translation" (a new handoff went out), the
# handoff/generation timestamp advances, and an error detail explaining the
# stale handback is recorded on the language sheets.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/caa38f1f6560b68b85d43ab0b59c7111648fc132/e2e/47407ed6-37c6-4239-8f04-d64dee116c62.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5705e66fcf870b4a46682795b0b6719418c7c935/e2e/47407ed6-37c6-4239-8f04-d64dee116c62.md."

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"
$wsOverview.Range("G3").Value = "2017-02-17 07:51:09"

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = "In Translation"
$wsZh.Range("H3").Value = "2017-02-17 07:50:51"
$wsZh.Range("R3").Value = $errorDetail

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = "In Translation"
$wsDe.Range("H3").Value = "2017-02-17 07:51:09"
$wsDe.Range("R3").Value = $errorDetail
